# Apply the update described in the commit:
#  - Extend the 72-month rolling temperature statistics range from row 75
#    to row 135 (i.e. H3:I75 -> H3:I135) by filling in new H (actual) and
#    I (average) values for rows 76 through 135.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics calculator")

# New "H" column (actual temperature) values for rows 76..135, in order.
$hValues = @(
    3.7530060000000001,6.3609549999999997,5.7115340000000003,10.877504,12.746848,
    15.079478999999999,16.463688000000001,18.133773999999999,13.475885999999999,8.9190860000000001,
    7.6639989999999996,1.808792,2.705409,4.0119020000000001,5.369224,
    6.1563179999999997,11.897114,15.144710999999999,18.398571,20.152441,
    15.909276999999999,9.7026559999999993,5.4410220000000002,4.0441940000000001,5.1618069999999996,
    3.9569990000000002,4.7896400000000003,7.6771289999999999,13.463158,14.901141000000001,
    19.623201000000002,18.838439999999999,13.972019,11.201335,6.6556240000000004,
    3.397316,4.5746159999999998,2.3380610000000002,4.606725,8.2302590000000002,
    12.049116,15.336563999999999,17.199770000000001,18.466093000000001,13.688967999999999,
    7.4050089999999997,7.978154,4.8641350000000001,3.8055110000000001,3.4569990000000002,
    4.6081190000000003,8.8570279999999997,11.587792,14.130141,18.377230000000001,
    18.955072000000001,17.859881999999999,11.608667000000001,5.2092049999999999,4.3596599999999999
)

# New "I" column (average reference) value - constant across the new rows.
$iValue = 7.8185159999999998

$startRow = 76
for ($idx = 0; $idx -lt $hValues.Length; $idx++) {
    $row = $startRow + $idx
    $ws.Cells.Item($row, 8).Value = $hValues[$idx]   # column H
    $ws.Cells.Item($row, 9).Value = $iValue           # column I
}

$endRow = $startRow + $hValues.Length - 1

# Update the selected range to reflect the extended data (H3:I75 -> H3:I135)
$ws.Range("H3:I$endRow").Select()
